$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.350.30"
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("D3").Value = "3.415.03"
$ws.Range("E3").Value = "  -1.49%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.24"
$ws.Range("E5").Value = "  -1.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.74"
$ws.Range("E6").Value = "  -3.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.610"
$ws.Range("E7").Value = "  +3.70%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "3.427.82"
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.20"
$ws.Range("E10").Value = "  -2.00%  "
$ws.Range("E11").Value = "  -3.64%  "
$ws.Range("E12").Value = "  -1.39%  "
$ws.Range("D13").Value = "4.003.38"
$ws.Range("E13").Value = "  -1.32%  "
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("E15").Value = "  -4.06%  "
$ws.Range("E16").Value = "  -4.07%  "
$ws.Range("D17").Value = "64.426.82"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "3.408.99"
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.34"
$ws.Range("E19").Value = "  -1.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.97"
$ws.Range("E20").Value = "  -4.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "374.12"
$ws.Range("E21").Value = "  -4.78%  "
$ws.Range("E22").Value = "  -3.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.549"
$ws.Range("E23").Value = "  -0.49%  "
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.14"
$ws.Range("E25").Value = "  -1.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000119"
$ws.Range("E26").Value = "  -5.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.03"
$ws.Range("E27").Value = "  +4.58%  "
$ws.Range("E28").Value = "  -2.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("E30").Value = "  +2.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.10"
$ws.Range("E31").Value = "  -1.88%  "
$ws.Range("E32").Value = "  -1.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.11"
$ws.Range("E33").Value = "  -2.88%  "
$ws.Range("E34").Value = "  +0.86%  "
$ws.Range("E35").Value = "  +5.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "160.28"
$ws.Range("E36").Value = "  -0.63%  "
$ws.Range("E37").Value = "  -0.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0757"
$ws.Range("E38").Value = "  -2.85%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.80"
$ws.Range("E39").Value = "  +3.12%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.74"
$ws.Range("E40").Value = "  -3.32%  "
$ws.Range("D41").Value = "2.852.94"
$ws.Range("E41").Value = "  -2.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.60"
$ws.Range("E42").Value = "  +0.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.67"
$ws.Range("E43").Value = "  -0.30%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "26.05"
$ws.Range("E44").Value = "  +8.35%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0313"
$ws.Range("E45").Value = "  -2.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.768"
$ws.Range("E46").Value = "  -1.10%  "
$ws.Range("E47").Value = "  -2.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "314.33"
$ws.Range("E48").Value = "  +5.47%  "
$ws.Range("E49").Value = "  +0.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.55"
$ws.Range("E50").Value = "  -1.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.851"
$ws.Range("E51").Value = "  -2.87%  "
